# Fix IFRS financial figures for 동양네트웍스ㅢ company_list sheet.
# The previously-entered values were off by orders of magnitude / wrong
# altogether; this replaces rows 2-9 with the corrected figures, drops the
# stray V3 value, and clears out rows 7-9 (which should only keep their
# leading A/B/C identifier columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (2014/12) ---
$ws.Range("D2").Value = 1430
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 1689
$ws.Range("H2").Value = 1435
$ws.Range("I2").Value = 1433
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 1278
$ws.Range("L2").Value = 925
$ws.Range("M2").Value = 353
$ws.Range("N2").Value = 350
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 157
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 620
$ws.Range("S2").Value = -524
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 11
$ws.Range("W2").Value = 0.53
$ws.Range("X2").Value = 100.32
$ws.Range("Y2").Value = -298.27
$ws.Range("Z2").Value = 81.61
$ws.Range("AA2").Value = 261.79
$ws.Range("AB2").Value = 92.16
$ws.Range("AC2").Value = 4306
$ws.Range("AD2").Value = 0.13
$ws.Range("AE2").Value = 868
$ws.Range("AF2").Value = 0.65
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 40357767

# --- Row 3 (2015/12) ---
$ws.Range("D3").Value = 1123
$ws.Range("E3").Value = -65
$ws.Range("F3").Value = -65
$ws.Range("G3").Value = -120
$ws.Range("H3").Value = -90
$ws.Range("I3").Value = -90
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 975
$ws.Range("L3").Value = 478
$ws.Range("M3").Value = 497
$ws.Range("N3").Value = 494
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 263
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 361
$ws.Range("S3").Value = -201
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = -3
$ws.Range("V3").ClearContents()
$ws.Range("W3").Value = -5.78
$ws.Range("X3").Value = -7.98
$ws.Range("Y3").Value = -21.25
$ws.Range("Z3").Value = -7.95
$ws.Range("AA3").Value = 96.23
$ws.Range("AB3").Value = 50.3
$ws.Range("AC3").Value = -220
$ws.Range("AD3").Value = -6.21
$ws.Range("AE3").Value = 938
$ws.Range("AF3").Value = 1.46
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 52636187

# --- Row 4 (2016/12) ---
$ws.Range("D4").Value = 941
$ws.Range("E4").Value = -66
$ws.Range("F4").Value = -66
$ws.Range("G4").Value = -21
$ws.Range("H4").Value = -40
$ws.Range("I4").Value = -40
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 964
$ws.Range("L4").Value = 548
$ws.Range("M4").Value = 416
$ws.Range("N4").Value = 407
$ws.Range("O4").Value = 9
$ws.Range("P4").Value = 264
$ws.Range("Q4").Value = -24
$ws.Range("R4").Value = 24
$ws.Range("S4").Value = 84
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = -28
$ws.Range("V4").Value = 218
$ws.Range("W4").Value = -7.06
$ws.Range("X4").Value = -4.22
$ws.Range("Y4").Value = -8.98
$ws.Range("Z4").Value = -4.09
$ws.Range("AA4").Value = 131.8
$ws.Range("AB4").Value = 34.18
$ws.Range("AC4").Value = -77
$ws.Range("AD4").Value = -14.54
$ws.Range("AE4").Value = 772
$ws.Range("AF4").Value = 1.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 52734346

# --- Row 5 (2017/12) ---
$ws.Range("D5").Value = 863
$ws.Range("E5").Value = -78
$ws.Range("F5").Value = -78
$ws.Range("G5").Value = -200
$ws.Range("H5").Value = -203
$ws.Range("I5").Value = -201
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 1018
$ws.Range("L5").Value = 230
$ws.Range("M5").Value = 788
$ws.Range("N5").Value = 781
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 455
$ws.Range("Q5").Value = -66
$ws.Range("R5").Value = 18
$ws.Range("S5").Value = 185
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = -67
$ws.Range("V5").Value = 27
$ws.Range("W5").Value = -9.06
$ws.Range("X5").Value = -23.58
$ws.Range("Y5").Value = -33.91
$ws.Range("Z5").Value = -20.52
$ws.Range("AA5").Value = 29.21
$ws.Range("AB5").Value = 60.34
$ws.Range("AC5").Value = -328
$ws.Range("AD5").Value = -4.81
$ws.Range("AE5").Value = 857
$ws.Range("AF5").Value = 1.84
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 91096176

# --- Row 6 (2018/12) ---
$ws.Range("D6").Value = 700
$ws.Range("E6").Value = -198
$ws.Range("F6").Value = -198
$ws.Range("G6").Value = -116
$ws.Range("H6").Value = -116
$ws.Range("I6").Value = -114
$ws.Range("K6").Value = 1906
$ws.Range("L6").Value = 1022
$ws.Range("M6").Value = 884
$ws.Range("N6").Value = 884
$ws.Range("P6").Value = 472
$ws.Range("Q6").Value = -283
$ws.Range("R6").Value = -557
$ws.Range("S6").Value = 962
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = -286
$ws.Range("V6").Value = 723
$ws.Range("W6").Value = -28.31
$ws.Range("X6").Value = -16.6
$ws.Range("Y6").Value = -13.74
$ws.Range("Z6").Value = -7.95
$ws.Range("AA6").Value = 115.6
$ws.Range("AB6").Value = 87.61
$ws.Range("AC6").Value = -121
$ws.Range("AD6").Value = -18.33
$ws.Range("AE6").Value = 937
$ws.Range("AF6").Value = 2.37
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 94336452

# --- Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) ---
# These forward estimate rows were erroneous altogether; clear every data
# column, keeping only the leading identifier columns A-C.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
